$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Marking row (row 11): Right = 9, Wrong = 2
$ws.Range("B11").Value = 9
$ws.Range("C11").Value = 2

# Update Total row (row 12): Right total = 117, and Max label "117/252"
$ws.Range("B12").Value = 117
$ws.Range("E12").Value = "117/252"
